# Insert 3 new price-record rows (Candy White / Especial, Primera, Segunda)
# just before the current row 336, pushing the existing rows 336-439 down to
# 339-442 (new dimension becomes A1:T442).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 336..439 down by 3 rows; the new blank rows 336-338
# inherit the formatting (incl. the date-format style) of the row above.
$ws.Rows("336:338").Insert()

# --- Row 336: Candy White / Especial ---
$ws.Cells.Item(336, 1).Value = 11
$ws.Cells.Item(336, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(336, 3).Value = "Bíobío"
$ws.Cells.Item(336, 4).Value = 44951
$ws.Cells.Item(336, 5).Value = 8
$ws.Cells.Item(336, 6).Value = "Fruta"
$ws.Cells.Item(336, 7).Value = 100103
$ws.Cells.Item(336, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(336, 9).Value = 100103006
$ws.Cells.Item(336, 10).Value = "Nectarín"
$ws.Cells.Item(336, 11).Value = "Candy White"
$ws.Cells.Item(336, 12).Value = "Especial"
$ws.Cells.Item(336, 13).Value = 120
$ws.Cells.Item(336, 14).Value = 14000
$ws.Cells.Item(336, 15).Value = 14000
$ws.Cells.Item(336, 16).Value = 14000
$ws.Cells.Item(336, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(336, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(336, 19).Value = 933
$ws.Cells.Item(336, 20).Value = 15

# --- Row 337: Candy White / Primera ---
$ws.Cells.Item(337, 1).Value = 11
$ws.Cells.Item(337, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(337, 3).Value = "Bíobío"
$ws.Cells.Item(337, 4).Value = 44951
$ws.Cells.Item(337, 5).Value = 8
$ws.Cells.Item(337, 6).Value = "Fruta"
$ws.Cells.Item(337, 7).Value = 100103
$ws.Cells.Item(337, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(337, 9).Value = 100103006
$ws.Cells.Item(337, 10).Value = "Nectarín"
$ws.Cells.Item(337, 11).Value = "Candy White"
$ws.Cells.Item(337, 12).Value = "Primera"
$ws.Cells.Item(337, 13).Value = 150
$ws.Cells.Item(337, 14).Value = 12000
$ws.Cells.Item(337, 15).Value = 12000
$ws.Cells.Item(337, 16).Value = 12000
$ws.Cells.Item(337, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(337, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(337, 19).Value = 800
$ws.Cells.Item(337, 20).Value = 15

# --- Row 338: Candy White / Segunda ---
$ws.Cells.Item(338, 1).Value = 11
$ws.Cells.Item(338, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(338, 3).Value = "Bíobío"
$ws.Cells.Item(338, 4).Value = 44951
$ws.Cells.Item(338, 5).Value = 8
$ws.Cells.Item(338, 6).Value = "Fruta"
$ws.Cells.Item(338, 7).Value = 100103
$ws.Cells.Item(338, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(338, 9).Value = 100103006
$ws.Cells.Item(338, 10).Value = "Nectarín"
$ws.Cells.Item(338, 11).Value = "Candy White"
$ws.Cells.Item(338, 12).Value = "Segunda"
$ws.Cells.Item(338, 13).Value = 100
$ws.Cells.Item(338, 14).Value = 10000
$ws.Cells.Item(338, 15).Value = 10000
$ws.Cells.Item(338, 16).Value = 10000
$ws.Cells.Item(338, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(338, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(338, 19).Value = 667
$ws.Cells.Item(338, 20).Value = 15
